# Regenerate merged AHB files
#
# 1. Rename the paired "_old" / "_new" comparison-column headers in row 1
#    to the concrete revision tags "_FV2304" / "_FV2310".
# 2. Turn the used range into a native Excel Table ("Table1") with an
#    AutoFilter, matching the renamed headers.
# 3. Freeze the header row (split/freeze at row 2, i.e. row 1 stays visible).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header renames -----------------------------------------------------
$headerMap = @{
    "A1" = "Segmentname_FV2304"
    "B1" = "Segmentgruppe_FV2304"
    "C1" = "Segment_FV2304"
    "D1" = "Datenelement_FV2304"
    "E1" = "Segment ID_FV2304"
    "F1" = "Code_FV2304"
    "G1" = "Qualifier_FV2304"
    "H1" = "Beschreibung_FV2304"
    "I1" = "Bedingungsausdruck_FV2304"
    "J1" = "Bedingung_FV2304"
    "L1" = "Segmentname_FV2310"
    "M1" = "Segmentgruppe_FV2310"
    "N1" = "Segment_FV2310"
    "O1" = "Datenelement_FV2310"
    "P1" = "Segment ID_FV2310"
    "Q1" = "Code_FV2310"
    "R1" = "Qualifier_FV2310"
    "S1" = "Beschreibung_FV2310"
    "T1" = "Bedingungsausdruck_FV2310"
    "U1" = "Bedingung_FV2310"
}

foreach ($addr in $headerMap.Keys) {
    $ws.Range($addr).Value = $headerMap[$addr]
}

# --- 2. Convert the used range A1:U72 into an Excel Table ------------------
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U72"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ----------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
